# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-45, columns C:E) is re-sorted:
# previously grouped by worker (5x NORVELIS, 7x KARINA, 6x MARLIZ,
# 7x NUBIS, 5x ANDREA) it becomes grouped by period (2302..2308), and two
# new workers (KARINA MARGARITA MONTES CARMONA / 1047455394 and NUBIS
# CAROLINA VERGARA SILGADO / 1001974820, etc. already existed) are woven
# into every period row so each of the 5 periodos 2304-2308 now lists all
# 5 trabajadores, while 2302/2303 list the subset that was already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, DocTrabajador (C), NombreTrabajador (D), PeriodoMora (E)
$rows = @(
    @(16, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2302'),
    @(17, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2302'),
    @(18, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2303'),
    @(19, '1047510112', 'MARLIZ ARRIETA JULIO', '2303'),
    @(20, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2303'),
    @(21, '1002244933', 'ANDREA MARCELA MARIMON CORREA', '2303'),
    @(22, '1148434315', 'NORVELIS DE ALBA ARRIETA', '2304'),
    @(23, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2304'),
    @(24, '1047510112', 'MARLIZ ARRIETA JULIO', '2304'),
    @(25, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2304'),
    @(26, '1002244933', 'ANDREA MARCELA MARIMON CORREA', '2304'),
    @(27, '1148434315', 'NORVELIS DE ALBA ARRIETA', '2305'),
    @(28, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2305'),
    @(29, '1047510112', 'MARLIZ ARRIETA JULIO', '2305'),
    @(30, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2305'),
    @(31, '1002244933', 'ANDREA MARCELA MARIMON CORREA', '2305'),
    @(32, '1148434315', 'NORVELIS DE ALBA ARRIETA', '2306'),
    @(33, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2306'),
    @(34, '1047510112', 'MARLIZ ARRIETA JULIO', '2306'),
    @(35, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2306'),
    @(36, '1002244933', 'ANDREA MARCELA MARIMON CORREA', '2306'),
    @(37, '1148434315', 'NORVELIS DE ALBA ARRIETA', '2307'),
    @(38, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2307'),
    @(39, '1047510112', 'MARLIZ ARRIETA JULIO', '2307'),
    @(40, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2307'),
    @(41, '1148434315', 'NORVELIS DE ALBA ARRIETA', '2308'),
    @(42, '1047455394', 'KARINA MARGARITA MONTES CARMONA', '2308'),
    @(43, '1047510112', 'MARLIZ ARRIETA JULIO', '2308'),
    @(44, '1001974820', 'NUBIS CAROLINA VERGARA SILGADO', '2308'),
    @(45, '1002244933', 'ANDREA MARCELA MARIMON CORREA', '2308')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
